$wb = $excel.ActiveWorkbook

# Both the "Summary" sheet and the "Pattern1-Pure Data" sheet carry the
# same row (row 3 = model "gemini-3-pro" under "Pattern1-Pure Data")
# that needs its metrics refreshed for the new data date.
$sheetNames = @("Summary", "Pattern1-Pure Data")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Plain currency-style text values are safe to assign directly since
    # the leading currency sign prevents Excel from reinterpreting them
    # as numbers.
    $ws.Range("D3").Value = "¥1,001,002.00"
    $ws.Range("E3").Value = "¥+1,002.00"

    # Percentage-looking text needs a leading apostrophe so Excel keeps
    # it as literal text instead of auto-converting it to a formatted
    # percentage number; ClearFormats() then strips the quote-prefix
    # style Excel applies so the cell keeps its original (unstyled)
    # appearance.
    $ws.Range("F3").Value = "'+0.10%"
    $ws.Range("F3").ClearFormats()

    $ws.Range("G3").Value = "'+28.71%"
    $ws.Range("G3").ClearFormats()

    $ws.Range("H3").Value = 0

    $ws.Range("I3").Value = "'0.00%"
    $ws.Range("I3").ClearFormats()

    $ws.Range("J3").Value = "'100.0%"
    $ws.Range("J3").ClearFormats()

    $ws.Range("K3").Value = "'0.1002%"
    $ws.Range("K3").ClearFormats()

    $ws.Range("L3").Value = "'0.0000%"
    $ws.Range("L3").ClearFormats()

    $ws.Range("M3").Value = 2

    # Numeric-looking date-code text also needs the apostrophe trick so
    # it stays a text string (e.g. "20251218") rather than becoming a
    # number.
    $ws.Range("O3").Value = "'20251218"
    $ws.Range("O3").ClearFormats()
}
